$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data to match the latest GitHub Actions scrape.
# D (Price) and E (Volume(1h)) columns store text-like values (e.g. "1.010", "  +0.26%  ")
# so each target cell is forced to Text format before assignment to avoid Excel
# auto-converting them to numbers/percentages.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.396.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.693.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.11%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.36%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.07"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5491"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.18%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2744"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.34%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06461"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.48%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.27%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07677"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.54%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.693.28"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.53%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.57%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5834"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.47%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008351"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.97%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.48"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.36%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.439.46"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.934"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.36%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.91"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.11%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.251"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.71"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.71%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1329"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.77%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.913"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.79"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06287"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.11%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.88%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.332"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.604"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.611"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.83%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.683"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.71%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.043"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.04%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.40%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.412"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.69%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01642"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.23%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.179"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.96%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.118.40"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.32%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8841"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.49%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.18%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.72"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.75%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.845.87"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.37%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.54"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000108"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.013"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.203"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05281"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.18%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4306"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.03%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.105"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.46%  "
